$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking values
# like "526.38" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.050.46"
$ws.Range("E2").Value = "  +3.02%  "

$ws.Range("D3").Value = "2.542.30"
$ws.Range("E3").Value = "  +5.30%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "526.38"
$ws.Range("E5").Value = "  +2.96%  "

$ws.Range("D6").Value = "135.22"
$ws.Range("E6").Value = "  +4.84%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +3.93%  "

$ws.Range("D9").Value = "2.540.93"
$ws.Range("E9").Value = "  +4.98%  "

$ws.Range("D10").Value = "0.0992"
$ws.Range("E10").Value = "  +3.95%  "

$ws.Range("E11").Value = "  -0.77%  "

$ws.Range("E12").Value = "  +1.03%  "

$ws.Range("E13").Value = "  +1.69%  "

$ws.Range("D14").Value = "2.996.12"
$ws.Range("E14").Value = "  +5.42%  "

$ws.Range("D15").Value = "59.153.95"
$ws.Range("E15").Value = "  +3.33%  "

$ws.Range("D16").Value = "22.42"
$ws.Range("E16").Value = "  +4.95%  "

$ws.Range("E17").Value = "  +3.88%  "

$ws.Range("D18").Value = "2.545.66"
$ws.Range("E18").Value = "  +5.34%  "

$ws.Range("D19").Value = "10.76"
$ws.Range("E19").Value = "  +4.26%  "

$ws.Range("D20").Value = "324.18"
$ws.Range("E20").Value = "  +3.54%  "

$ws.Range("E21").Value = "  +3.51%  "

$ws.Range("D22").Value = "6.11"
$ws.Range("E22").Value = "  +8.84%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").Value = "65.32"
$ws.Range("E24").Value = "  +2.97%  "

$ws.Range("E25").Value = "  +2.82%  "

$ws.Range("E26").Value = "  +0.31%  "

$ws.Range("E27").Value = "  +1.58%  "

$ws.Range("E28").Value = "  +5.17%  "

$ws.Range("D29").Value = "0.0₃0761"
$ws.Range("E29").Value = "  +6.18%  "

$ws.Range("E30").Value = "  +7.85%  "

$ws.Range("E31").Value = "  +5.10%  "

$ws.Range("D32").Value = "169.51"
$ws.Range("E32").Value = "  +0.23%  "

$ws.Range("E33").Value = "  +3.59%  "

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").Value = "18.28"
$ws.Range("E36").Value = "  +3.71%  "

$ws.Range("E37").Value = "  +0.61%  "

$ws.Range("E38").Value = "  +3.89%  "

$ws.Range("E39").Value = "  +6.17%  "

$ws.Range("D40").Value = "36.76"
$ws.Range("E40").Value = "  +1.39%  "

$ws.Range("E41").Value = "  +3.34%  "

$ws.Range("E42").Value = "  +6.26%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "134.29"
$ws.Range("E43").Value = "  +10.30%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "3.49"
$ws.Range("E44").Value = "  +4.25%  "

$ws.Range("D45").Value = "5.12"
$ws.Range("E45").Value = "  +5.39%  "

$ws.Range("D46").Value = "0.604"
$ws.Range("E46").Value = "  +4.03%  "

$ws.Range("E47").Value = "  +2.96%  "

$ws.Range("D48").Value = "0.0507"
$ws.Range("E48").Value = "  +6.00%  "

$ws.Range("E49").Value = "  +4.29%  "

$ws.Range("D50").Value = "17.17"

$ws.Range("D51").Value = "1.759.69"
$ws.Range("E51").Value = "  +4.47%  "

# Restore default styling so no stray formatting is introduced.
$ws.Range("D2:D51").Style = "Normal"
